$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1372.5745
$ws.Range("I15").Value = 1372.5745
$ws.Range("K15").Value = 4117.7235
$ws.Range("M15").Value = -3948.7235
$ws.Range("H32").Value = 1359.875
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H33").Value = 120.7
$ws.Range("I33").Value = 121.75
$ws.Range("J33").Value = 120
$ws.Range("K33").Value = 121.75
$ws.Range("L33").Value = 120
$ws.Range("M33").Value = 107.25
$ws.Range("N33").Value = -578
$ws.Range("H51").Value = 6000
$ws.Range("I51").Value = 5000
$ws.Range("K51").Value = 5000
$ws.Range("M51").Value = -4516
$ws.Range("H98").Value = 2309.476
$ws.Range("I98").Value = 2447.3684
$ws.Range("J98").Value = 999.5
$ws.Range("K98").Value = 2447.3684
$ws.Range("L98").Value = 999.5
$ws.Range("M98").Value = -949.3683999999998
$ws.Range("N98").Value = -3995.5
$ws.Range("H121").Value = 1075
$ws.Range("J121").Value = 2000
$ws.Range("L121").Value = 6000
$ws.Range("N121").Value = -9494
$ws.Range("H122").Value = 2309.476
$ws.Range("I122").Value = 2447.3684
$ws.Range("J122").Value = 999.5
$ws.Range("K122").Value = 7342.1052
$ws.Range("L122").Value = 2998.5
$ws.Range("M122").Value = -4892.1052
$ws.Range("N122").Value = -7898.5
$ws.Range("H125").Value = 486.33334
$ws.Range("I125").Value = 383.6
$ws.Range("J125").Value = 1000
$ws.Range("K125").Value = 3452.4
$ws.Range("L125").Value = 9000
$ws.Range("M125").Value = -992.4000000000001
$ws.Range("N125").Value = -13920
$ws.Range("H137").Value = 1457.8
$ws.Range("I137").Value = 1402.3
$ws.Range("J137").Value = 1679.8
$ws.Range("K137").Value = 4206.9
$ws.Range("L137").Value = 5039.4
$ws.Range("M137").Value = -1656.9
$ws.Range("N137").Value = -10139.4
$ws.Range("H138").Value = 3166.7073
$ws.Range("J138").Value = 4023.6
$ws.Range("L138").Value = 12070.8
$ws.Range("N138").Value = -22350.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 723.8
$ws.Range("I2").Value = 624.8
$ws.Range("J2").Value = 921.8
$ws.Range("K2").Value = 624.8
$ws.Range("L2").Value = 921.8
$ws.Range("M2").Value = -511.8
$ws.Range("N2").Value = -1147.8
$ws.Range("H45").Value = 1750.7222
$ws.Range("I45").Value = 1647
$ws.Range("K45").Value = 1647
$ws.Range("M45").Value = -1270
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H116").Value = 723.8
$ws.Range("I116").Value = 624.8
$ws.Range("J116").Value = 921.8
$ws.Range("K116").Value = 624.8
$ws.Range("L116").Value = 921.8
$ws.Range("M116").Value = 1669.2
$ws.Range("N116").Value = -5509.8
$ws.Range("H132").Value = 2088.8462
$ws.Range("I132").Value = 1230.375
$ws.Range("J132").Value = 3462.4
$ws.Range("K132").Value = 3691.125
$ws.Range("L132").Value = 10387.2
$ws.Range("M132").Value = -1161.125
$ws.Range("N132").Value = -15447.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 723.8
$ws.Range("I3").Value = 624.8
$ws.Range("J3").Value = 921.8
$ws.Range("K3").Value = 624.8
$ws.Range("L3").Value = 921.8
$ws.Range("M3").Value = -510.8
$ws.Range("N3").Value = -1149.8
$ws.Range("H17").Value = 57009
$ws.Range("J17").Value = 57009
$ws.Range("L17").Value = 57009
$ws.Range("N17").Value = -57353
$ws.Range("H105").Value = 1985.6818
$ws.Range("I105").Value = 2019.6471
$ws.Range("K105").Value = 2019.6471
$ws.Range("M105").Value = -272.6470999999999
$ws.Range("H107").Value = 3524.4
$ws.Range("I107").Value = 3524.4
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3524.4
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1604.4
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H22").Value = 699.8
$ws.Range("I22").Value = 374.75
$ws.Range("K22").Value = 374.75
$ws.Range("M22").Value = -24.75
$ws.Range("H31").Value = 3617.4285
$ws.Range("I31").Value = 2752.375
$ws.Range("J31").Value = 4770.8335
$ws.Range("K31").Value = 2752.375
$ws.Range("L31").Value = 4770.8335
$ws.Range("M31").Value = -2457.375
$ws.Range("N31").Value = -5360.8335
$ws.Range("H34").Value = 3617.4285
$ws.Range("I34").Value = 2752.375
$ws.Range("J34").Value = 4770.8335
$ws.Range("K34").Value = 2752.375
$ws.Range("L34").Value = 4770.8335
$ws.Range("M34").Value = -2550.375
$ws.Range("N34").Value = -5174.8335
$ws.Range("H58").Value = 1788.0769
$ws.Range("I58").Value = 1479.25
$ws.Range("J58").Value = 2282.2
$ws.Range("K58").Value = 1479.25
$ws.Range("L58").Value = 2282.2
$ws.Range("M58").Value = -1276.25
$ws.Range("N58").Value = -2688.2
$ws.Range("H107").Value = 311.36365
$ws.Range("I107").Value = 255.26315
$ws.Range("J107").Value = 666.6667
$ws.Range("K107").Value = 255.26315
$ws.Range("L107").Value = 666.6667
$ws.Range("M107").Value = 1664.73685
$ws.Range("N107").Value = -4506.6667
$ws.Range("H134").Value = 1172.5454
$ws.Range("I134").Value = 985.55554
$ws.Range("J134").Value = 2014
$ws.Range("K134").Value = 2956.66662
$ws.Range("L134").Value = 6042
$ws.Range("M134").Value = -421.66662
$ws.Range("N134").Value = -11112
$ws.Range("H136").Value = 1788.0769
$ws.Range("I136").Value = 1479.25
$ws.Range("J136").Value = 2282.2
$ws.Range("K136").Value = 4437.75
$ws.Range("L136").Value = 6846.599999999999
$ws.Range("M136").Value = -1887.75
$ws.Range("N136").Value = -11946.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 102.23077
$ws.Range("I2").Value = 98.09090999999999
$ws.Range("J2").Value = 125
$ws.Range("K2").Value = 588.5454599999999
$ws.Range("L2").Value = 750
$ws.Range("M2").Value = -475.5454599999999
$ws.Range("N2").Value = -976
$ws.Range("H6").Value = 131.75
$ws.Range("I6").Value = 42.333332
$ws.Range("J6").Value = 400
$ws.Range("K6").Value = 126.999996
$ws.Range("L6").Value = 1200
$ws.Range("M6").Value = -13.999996
$ws.Range("N6").Value = -1426
$ws.Range("H33").Value = 111.333336
$ws.Range("J33").Value = 124.833336
$ws.Range("L33").Value = 749.000016
$ws.Range("N33").Value = -1315.000016
$ws.Range("H113").Value = 12010.111
$ws.Range("I113").Value = 33867
$ws.Range("K113").Value = 101601
$ws.Range("M113").Value = -99431
$ws.Range("H131").Value = 747.36
$ws.Range("J131").Value = 774.0879
$ws.Range("L131").Value = 2322.2637
$ws.Range("N131").Value = -12402.2637

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3996.1428
$ws.Range("I102").Value = 4162
$ws.Range("K102").Value = 4162
$ws.Range("M102").Value = -2540
$ws.Range("H107").Value = 932.5
$ws.Range("J107").Value = 1099
$ws.Range("L107").Value = 1099
$ws.Range("N107").Value = -4939
$ws.Range("H126").Value = 35932.965
$ws.Range("I126").Value = 2826.2856
$ws.Range("J126").Value = 113181.89
$ws.Range("K126").Value = 8478.856800000001
$ws.Range("L126").Value = 339545.67
$ws.Range("M126").Value = -6008.856800000001
$ws.Range("N126").Value = -344485.67
$ws.Range("H132").Value = 3447.1667
$ws.Range("I132").Value = 3447.1667
$ws.Range("K132").Value = 10341.5001
$ws.Range("M132").Value = -7811.500100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 12498.333
$ws.Range("I122").Value = 12080.833
$ws.Range("J122").Value = 13333.333
$ws.Range("K122").Value = 36242.499
$ws.Range("L122").Value = 39999.999
$ws.Range("M122").Value = -33792.499
$ws.Range("N122").Value = -44899.999
$ws.Range("H136").Value = 3657.4375
$ws.Range("I136").Value = 2002.5
$ws.Range("J136").Value = 5312.375
$ws.Range("K136").Value = 6007.5
$ws.Range("L136").Value = 15937.125
$ws.Range("M136").Value = -3457.5
$ws.Range("N136").Value = -21037.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1518.6666
$ws.Range("I81").Value = 1304.8889
$ws.Range("J81").Value = 2160
$ws.Range("K81").Value = 2609.7778
$ws.Range("L81").Value = 4320
$ws.Range("M81").Value = -1548.7778
$ws.Range("N81").Value = -6442
$ws.Range("H84").Value = 1518.6666
$ws.Range("I84").Value = 1304.8889
$ws.Range("J84").Value = 2160
$ws.Range("K84").Value = 13048.889
$ws.Range("L84").Value = 21600
$ws.Range("M84").Value = -7744.888999999999
$ws.Range("N84").Value = -32208
$ws.Range("H97").Value = 39995
$ws.Range("J97").Value = 39995
$ws.Range("L97").Value = 39995
$ws.Range("N97").Value = -41977
$ws.Range("H136").Value = 2986.7917
$ws.Range("I136").Value = 3513.0833
$ws.Range("J136").Value = 2460.5
$ws.Range("K136").Value = 10539.2499
$ws.Range("L136").Value = 7381.5
$ws.Range("M136").Value = -7989.249899999999
$ws.Range("N136").Value = -12481.5
